$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1/2: swap Nom / Prénom, refresh names with "1" suffix
$ws.Range("A1").Value = "Prénom"
$ws.Range("B1").Value = "Dimitri1"

$ws.Range("A2").Value = "Nom"
$ws.Range("B2").Value = "Lefebvre1"

# Row 3: refreshed job title
$ws.Range("B3").Value = "Data Analyst1"

# Row 5/6: refreshed links (set before the catch phrase so the shared
# string table ends up in the same order as produced by Excel)
$ws.Range("B5").Value = "https://github.com/Dim29601"
$ws.Range("B6").Value = "https://www.linkedin.com/in/dim-lefebvre601"

# Row 4: refreshed catch phrase
$ws.Range("B4").Value = "1J’allie expertise technique et vision stratégique pour transformer vos données en performance. Mon expérience passée fait la différence."

# Update the active selection to mirror the saved view state
$ws.Range("B15").Select()
